{"js": "// Replace the inline \"Gross Plot Ratio for sites where vesting is required\"\n// picture with a plain hyperlink run (same text as the link target), leaving\n// the rest of the document (including the gross-plot-ratio bookmark that\n// wraps the whole body) untouched.\n\nconst url =\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F01_Gross_Plot_Ratio.jpg?h=100%25&w=100%25\";\n\nconst body = context.document.body;\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length === 0) {\n  throw new Error(\"No inline picture found to convert into a hyperlink.\");\n}\n\n// The doc only has the one picture, but be defensive and pick it out by its\n// alt-text description in case more pictures are ever present.\nlet pic = pictures.items[0];\nfor (let i = 0; i < pictures.items.length; i++) {\n  const candidate = pictures.items[i];\n  candidate.load(\"altTextDescription\");\n}\nawait context.sync();\nfor (let i = 0; i < pictures.items.length; i++) {\n  if (\n    pictures.items[i].altTextDescription ===\n    \"Gross Plot Ratio for sites where vesting is required\"\n  ) {\n    pic = pictures.items[i];\n    break;\n  }\n}\n\n// The picture is its own run inside a BodyText paragraph; grab the\n// paragraph so we can drop the link text into the same (now-empty) spot.\nconst paragraph = pic.paragraph;\n\n// Remove the picture run.\npic.delete();\nawait context.sync();\n\n// Insert the URL as the paragraph's text, then turn that text into a\n// hyperlink \u2014 setting Range.hyperlink both wraps the run in a\n// <w:hyperlink> that targets the URL and applies the built-in \"Hyperlink\"\n// character style to the run, matching Word's normal \"insert hyperlink\"\n// behavior.\nparagraph.insertText(url, Word.InsertLocation.start);\nawait context.sync();\n\nconst linkRange = paragraph.search(url, { matchCase: true }).getFirstOrNullObject();\nawait context.sync();\nlinkRange.hyperlink = url;\nawait context.sync();\n", "ps1": "# Replace the inline \"Gross Plot Ratio for sites where vesting is required\"\n# picture with a plain hyperlink run (display text == link target), leaving\n# the rest of the document (including the gross-plot-ratio bookmark that\n# wraps the whole body) untouched.\n\n$d = $word.ActiveDocument\n$url = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F01_Gross_Plot_Ratio.jpg?h=100%25&w=100%25\"\n\n# Locate the inline picture by its alt text (there is only one in this\n# document, but match on the description to be defensive/explicit).\n$target = $null\nforeach ($shape in $d.InlineShapes) {\n    if ($shape.AlternativeText -eq \"Gross Plot Ratio for sites where vesting is required\") {\n        $target = $shape\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"No inline picture found to convert into a hyperlink.\"\n}\n\n# The picture occupies its own run inside a BodyText paragraph. Replacing\n# the shape's Range text swaps the picture out for the link text in place,\n# then Hyperlinks.Add wraps that range in a <w:hyperlink> pointing at the\n# URL and applies the built-in \"Hyperlink\" character style to the run -\n# exactly what Word's Insert Hyperlink does.\n$range = $target.Range\n$range.Text = $url\n$d.Hyperlinks.Add($range, $url) | Out-Null\n"}
